$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(58, 8).Value = 1081.5  # H58: '765.55554' -> '1081.5'
$ws.Cells.Item(58, 9).Value = 901.6667  # I58: '765.55554' -> '901.6667'
$ws.Cells.Item(58, 10).Value = 1189.4  # J58: '0' -> '1189.4'
$ws.Cells.Item(58, 11).Value = 2705.0001  # K58: '2296.66662' -> '2705.0001'
$ws.Cells.Item(58, 12).Value = 3568.2  # L58: '0' -> '3568.2'
$ws.Cells.Item(58, 13).Value = -2555.0001  # M58: '-2146.66662' -> '-2555.0001'
$ws.Cells.Item(58, 14).Value = -3868.2  # N58: None -> '-3868.2'

$ws.Cells.Item(64, 8).Value = 3582.3044  # H64: '3491.662' -> '3582.3044'
$ws.Cells.Item(64, 9).Value = 3391.4285  # I64: '3247.6' -> '3391.4285'
$ws.Cells.Item(64, 11).Value = 3391.4285  # K64: '3247.6' -> '3391.4285'
$ws.Cells.Item(64, 13).Value = -3143.4285  # M64: '-2999.6' -> '-3143.4285'

$ws.Cells.Item(67, 8).Value = 3582.3044  # H67: '3491.662' -> '3582.3044'
$ws.Cells.Item(67, 9).Value = 3391.4285  # I67: '3247.6' -> '3391.4285'
$ws.Cells.Item(67, 11).Value = 3391.4285  # K67: '3247.6' -> '3391.4285'
$ws.Cells.Item(67, 13).Value = -2533.4285  # M67: '-2389.6' -> '-2533.4285'

$ws.Cells.Item(74, 8).Value = 3202.9312  # H74: '3339.5862' -> '3202.9312'
$ws.Cells.Item(74, 9).Value = 2397  # I74: '2737.6' -> '2397'
$ws.Cells.Item(74, 10).Value = 3370.8333  # J74: '3465' -> '3370.8333'
$ws.Cells.Item(74, 11).Value = 2397  # K74: '2737.6' -> '2397'
$ws.Cells.Item(74, 12).Value = 3370.8333  # L74: '3465' -> '3370.8333'
$ws.Cells.Item(74, 13).Value = -1461  # M74: '-1801.6' -> '-1461'
$ws.Cells.Item(74, 14).Value = -5242.8333  # N74: '-5337' -> '-5242.8333'

$ws.Cells.Item(77, 8).Value = 3202.9312  # H77: '3339.5862' -> '3202.9312'
$ws.Cells.Item(77, 9).Value = 2397  # I77: '2737.6' -> '2397'
$ws.Cells.Item(77, 10).Value = 3370.8333  # J77: '3465' -> '3370.8333'
$ws.Cells.Item(77, 11).Value = 11985  # K77: '13688' -> '11985'
$ws.Cells.Item(77, 12).Value = 16854.1665  # L77: '17325' -> '16854.1665'
$ws.Cells.Item(77, 13).Value = -7305  # M77: '-9008' -> '-7305'
$ws.Cells.Item(77, 14).Value = -26214.1665  # N77: '-26685' -> '-26214.1665'

$ws.Cells.Item(93, 8).Value = 83096.37  # H93: '79204' -> '83096.37'
$ws.Cells.Item(93, 10).Value = 83096.37  # J93: '79204' -> '83096.37'
$ws.Cells.Item(93, 12).Value = 83096.37  # L93: '79204' -> '83096.37'
$ws.Cells.Item(93, 14).Value = -88088.37  # N93: '-84196' -> '-88088.37'

$ws.Cells.Item(117, 8).Value = 58000  # H117: '24980' -> '58000'
$ws.Cells.Item(117, 10).Value = 58000  # J117: '24980' -> '58000'
$ws.Cells.Item(117, 12).Value = 58000  # L117: '24980' -> '58000'
$ws.Cells.Item(117, 14).Value = -67178  # N117: '-34158' -> '-67178'

$ws.Cells.Item(121, 8).Value = 930  # H121: '700.38464' -> '930'
$ws.Cells.Item(121, 10).Value = 930  # J121: '700.38464' -> '930'
$ws.Cells.Item(121, 12).Value = 2790  # L121: '2101.15392' -> '2790'
$ws.Cells.Item(121, 14).Value = -6284  # N121: '-5595.15392' -> '-6284'

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 1828.2667  # H61: '2024' -> '1828.2667'
$ws.Cells.Item(61, 9).Value = 1281.36  # I61: '1428.1' -> '1281.36'
$ws.Cells.Item(61, 10).Value = 4562.8  # J61: '5003.5' -> '4562.8'
$ws.Cells.Item(61, 11).Value = 1281.36  # K61: '1428.1' -> '1281.36'
$ws.Cells.Item(61, 12).Value = 4562.8  # L61: '5003.5' -> '4562.8'
$ws.Cells.Item(61, 13).Value = -1069.36  # M61: '-1216.1' -> '-1069.36'
$ws.Cells.Item(61, 14).Value = -4986.8  # N61: '-5427.5' -> '-4986.8'

$ws.Cells.Item(74, 8).Value = 1557.75  # H74: '1760.2858' -> '1557.75'
$ws.Cells.Item(74, 9).Value = 1117.7742  # I74: '1303.3334' -> '1117.7742'
$ws.Cells.Item(74, 10).Value = 4285.6  # J74: '4502' -> '4285.6'
$ws.Cells.Item(74, 11).Value = 1117.7742  # K74: '1303.3334' -> '1117.7742'
$ws.Cells.Item(74, 12).Value = 4285.6  # L74: '4502' -> '4285.6'
$ws.Cells.Item(74, 13).Value = -243.7742000000001  # M74: '-429.3334' -> '-243.7742000000001'
$ws.Cells.Item(74, 14).Value = -6033.6  # N74: '-6250' -> '-6033.6'

$ws.Cells.Item(77, 8).Value = 1557.75  # H77: '1760.2858' -> '1557.75'
$ws.Cells.Item(77, 9).Value = 1117.7742  # I77: '1303.3334' -> '1117.7742'
$ws.Cells.Item(77, 10).Value = 4285.6  # J77: '4502' -> '4285.6'
$ws.Cells.Item(77, 11).Value = 5588.871  # K77: '6516.666999999999' -> '5588.871'
$ws.Cells.Item(77, 12).Value = 21428  # L77: '22510' -> '21428'
$ws.Cells.Item(77, 13).Value = -1220.871  # M77: '-2148.666999999999' -> '-1220.871'
$ws.Cells.Item(77, 14).Value = -30164  # N77: '-31246' -> '-30164'

$ws.Cells.Item(119, 8).Value = 22349  # H119: '29849' -> '22349'
$ws.Cells.Item(119, 10).Value = 22349  # J119: '29849' -> '22349'
$ws.Cells.Item(119, 12).Value = 22349  # L119: '29849' -> '22349'
$ws.Cells.Item(119, 14).Value = -32025  # N119: '-39525' -> '-32025'

$ws.Cells.Item(124, 8).Value = 27569.5  # H124: '10548.846' -> '27569.5'
$ws.Cells.Item(124, 10).Value = 27569.5  # J124: '10548.846' -> '27569.5'
$ws.Cells.Item(124, 12).Value = 27569.5  # L124: '10548.846' -> '27569.5'
$ws.Cells.Item(124, 14).Value = -37389.5  # N124: '-20368.846' -> '-37389.5'

$ws.Cells.Item(125, 8).Value = 46359.168  # H125: '90143.336' -> '46359.168'
$ws.Cells.Item(125, 10).Value = 46359.168  # J125: '90143.336' -> '46359.168'
$ws.Cells.Item(125, 12).Value = 46359.168  # L125: '90143.336' -> '46359.168'
$ws.Cells.Item(125, 14).Value = -56199.168  # N125: '-99983.336' -> '-56199.168'

$ws.Cells.Item(132, 8).Value = 2039.4445  # H132: '1591.6538' -> '2039.4445'
$ws.Cells.Item(132, 9).Value = 1263.0714  # I132: '1075.9474' -> '1263.0714'
$ws.Cells.Item(132, 10).Value = 4756.75  # J132: '2991.4285' -> '4756.75'
$ws.Cells.Item(132, 11).Value = 3789.2142  # K132: '3227.8422' -> '3789.2142'
$ws.Cells.Item(132, 12).Value = 14270.25  # L132: '8974.2855' -> '14270.25'
$ws.Cells.Item(132, 13).Value = -1259.2142  # M132: '-697.8422' -> '-1259.2142'
$ws.Cells.Item(132, 14).Value = -19330.25  # N132: '-14034.2855' -> '-19330.25'

$ws.Cells.Item(136, 8).Value = 1828.2667  # H136: '2024' -> '1828.2667'
$ws.Cells.Item(136, 9).Value = 1281.36  # I136: '1428.1' -> '1281.36'
$ws.Cells.Item(136, 10).Value = 4562.8  # J136: '5003.5' -> '4562.8'
$ws.Cells.Item(136, 11).Value = 3844.08  # K136: '4284.299999999999' -> '3844.08'
$ws.Cells.Item(136, 12).Value = 13688.4  # L136: '15010.5' -> '13688.4'
$ws.Cells.Item(136, 13).Value = -1294.08  # M136: '-1734.299999999999' -> '-1294.08'
$ws.Cells.Item(136, 14).Value = -18788.4  # N136: '-20110.5' -> '-18788.4'

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 1891.3043  # H86: '1819.6' -> '1891.3043'
$ws.Cells.Item(86, 9).Value = 1724  # I86: '1681.579' -> '1724'
$ws.Cells.Item(86, 10).Value = 2686  # J86: '2058' -> '2686'
$ws.Cells.Item(86, 11).Value = 1724  # K86: '1681.579' -> '1724'
$ws.Cells.Item(86, 12).Value = 2686  # L86: '2058' -> '2686'
$ws.Cells.Item(86, 13).Value = -601  # M86: '-558.579' -> '-601'
$ws.Cells.Item(86, 14).Value = -4932  # N86: '-4304' -> '-4932'

$ws.Cells.Item(89, 8).Value = 1891.3043  # H89: '1819.6' -> '1891.3043'
$ws.Cells.Item(89, 9).Value = 1724  # I89: '1681.579' -> '1724'
$ws.Cells.Item(89, 10).Value = 2686  # J89: '2058' -> '2686'
$ws.Cells.Item(89, 11).Value = 8620  # K89: '8407.895' -> '8620'
$ws.Cells.Item(89, 12).Value = 13430  # L89: '10290' -> '13430'
$ws.Cells.Item(89, 13).Value = -3004  # M89: '-2791.895' -> '-3004'
$ws.Cells.Item(89, 14).Value = -24662  # N89: '-21522' -> '-24662'

$ws.Cells.Item(94, 8).Value = 502.05554  # H94: '565' -> '502.05554'
$ws.Cells.Item(94, 9).Value = 436.08334  # I94: '502.875' -> '436.08334'
$ws.Cells.Item(94, 10).Value = 634  # J94: '675.44446' -> '634'
$ws.Cells.Item(94, 11).Value = 436.08334  # K94: '502.875' -> '436.08334'
$ws.Cells.Item(94, 12).Value = 634  # L94: '675.44446' -> '634'
$ws.Cells.Item(94, 13).Value = 14.91665999999998  # M94: '-51.875' -> '14.91665999999998'
$ws.Cells.Item(94, 14).Value = -1536  # N94: '-1577.44446' -> '-1536'

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2599.6667  # H31: '2723.966' -> '2599.6667'
$ws.Cells.Item(31, 9).Value = 1656.1163  # I31: '1708.4878' -> '1656.1163'
$ws.Cells.Item(31, 10).Value = 4628.3  # J31: '5037' -> '4628.3'
$ws.Cells.Item(31, 11).Value = 1656.1163  # K31: '1708.4878' -> '1656.1163'
$ws.Cells.Item(31, 12).Value = 4628.3  # L31: '5037' -> '4628.3'
$ws.Cells.Item(31, 13).Value = -1361.1163  # M31: '-1413.4878' -> '-1361.1163'
$ws.Cells.Item(31, 14).Value = -5218.3  # N31: '-5627' -> '-5218.3'

$ws.Cells.Item(34, 8).Value = 2599.6667  # H34: '2723.966' -> '2599.6667'
$ws.Cells.Item(34, 9).Value = 1656.1163  # I34: '1708.4878' -> '1656.1163'
$ws.Cells.Item(34, 10).Value = 4628.3  # J34: '5037' -> '4628.3'
$ws.Cells.Item(34, 11).Value = 1656.1163  # K34: '1708.4878' -> '1656.1163'
$ws.Cells.Item(34, 12).Value = 4628.3  # L34: '5037' -> '4628.3'
$ws.Cells.Item(34, 13).Value = -1454.1163  # M34: '-1506.4878' -> '-1454.1163'
$ws.Cells.Item(34, 14).Value = -5032.3  # N34: '-5441' -> '-5032.3'

$ws.Cells.Item(63, 8).Value = 17500  # H63: '20000' -> '17500'
$ws.Cells.Item(63, 9).Value = 15000  # I63: '0' -> '15000'
$ws.Cells.Item(63, 11).Value = 15000  # K63: '0' -> '15000'
$ws.Cells.Item(63, 13).Value = -14314  # M63: None -> '-14314'

$ws.Cells.Item(66, 8).Value = 17500  # H66: '20000' -> '17500'
$ws.Cells.Item(66, 9).Value = 15000  # I66: '0' -> '15000'
$ws.Cells.Item(66, 11).Value = 45000  # K66: '0' -> '45000'
$ws.Cells.Item(66, 13).Value = -41568  # M66: None -> '-41568'

$ws.Cells.Item(132, 8).Value = 1676.6976  # H132: '1656.1818' -> '1676.6976'
$ws.Cells.Item(132, 9).Value = 1321.3667  # I132: '1322.0333' -> '1321.3667'
$ws.Cells.Item(132, 10).Value = 2496.6924  # J132: '2372.2144' -> '2496.6924'
$ws.Cells.Item(132, 11).Value = 3964.1001  # K132: '3966.0999' -> '3964.1001'
$ws.Cells.Item(132, 12).Value = 7490.0772  # L132: '7116.6432' -> '7490.0772'
$ws.Cells.Item(132, 13).Value = -1434.1001  # M132: '-1436.0999' -> '-1434.1001'
$ws.Cells.Item(132, 14).Value = -12550.0772  # N132: '-12176.6432' -> '-12550.0772'

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(62, 8).Value = 6757  # H62: '3553' -> '6757'
$ws.Cells.Item(62, 9).Value = 0  # I62: '499' -> '0'
$ws.Cells.Item(62, 10).Value = 6757  # J62: '6607' -> '6757'
$ws.Cells.Item(62, 11).Value = 0  # K62: '1497' -> '0'
$ws.Cells.Item(62, 12).Value = 20271  # L62: '19821' -> '20271'
$ws.Cells.Item(62, 13).ClearContents()  # M62: '-811' -> removed
$ws.Cells.Item(62, 14).Value = -21643  # N62: '-21193' -> '-21643'

$ws.Cells.Item(65, 8).Value = 6757  # H65: '3553' -> '6757'
$ws.Cells.Item(65, 9).Value = 0  # I65: '499' -> '0'
$ws.Cells.Item(65, 10).Value = 6757  # J65: '6607' -> '6757'
$ws.Cells.Item(65, 11).Value = 0  # K65: '4491' -> '0'
$ws.Cells.Item(65, 12).Value = 60813  # L65: '59463' -> '60813'
$ws.Cells.Item(65, 13).ClearContents()  # M65: '-1059' -> removed
$ws.Cells.Item(65, 14).Value = -67677  # N65: '-66327' -> '-67677'

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 8099.5293  # H70: '7241.4' -> '8099.5293'
$ws.Cells.Item(70, 9).Value = 9149.333000000001  # I70: '8377.5' -> '9149.333000000001'
$ws.Cells.Item(70, 10).Value = 5580  # J70: '4590.5' -> '5580'
$ws.Cells.Item(70, 11).Value = 9149.333000000001  # K70: '8377.5' -> '9149.333000000001'
$ws.Cells.Item(70, 12).Value = 5580  # L70: '4590.5' -> '5580'
$ws.Cells.Item(70, 13).Value = -8879.333000000001  # M70: '-8107.5' -> '-8879.333000000001'
$ws.Cells.Item(70, 14).Value = -6120  # N70: '-5130.5' -> '-6120'

$ws.Cells.Item(73, 8).Value = 8099.5293  # H73: '7241.4' -> '8099.5293'
$ws.Cells.Item(73, 9).Value = 9149.333000000001  # I73: '8377.5' -> '9149.333000000001'
$ws.Cells.Item(73, 10).Value = 5580  # J73: '4590.5' -> '5580'
$ws.Cells.Item(73, 11).Value = 9149.333000000001  # K73: '8377.5' -> '9149.333000000001'
$ws.Cells.Item(73, 12).Value = 5580  # L73: '4590.5' -> '5580'
$ws.Cells.Item(73, 13).Value = -8213.333000000001  # M73: '-7441.5' -> '-8213.333000000001'
$ws.Cells.Item(73, 14).Value = -7452  # N73: '-6462.5' -> '-7452'

$ws.Cells.Item(103, 8).Value = 0  # H103: '20000' -> '0'
$ws.Cells.Item(103, 10).Value = 0  # J103: '20000' -> '0'
$ws.Cells.Item(103, 12).Value = 0  # L103: '20000' -> '0'
$ws.Cells.Item(103, 14).ClearContents()  # N103: '-22344' -> removed

$ws.Cells.Item(132, 8).Value = 3909.8235  # H132: '3654.973' -> '3909.8235'
$ws.Cells.Item(132, 9).Value = 3662.0476  # I132: '3300.125' -> '3662.0476'
$ws.Cells.Item(132, 11).Value = 10986.1428  # K132: '9900.375' -> '10986.1428'
$ws.Cells.Item(132, 13).Value = -8456.1428  # M132: '-7370.375' -> '-8456.1428'

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 2628.1365  # H136: '2890.718' -> '2628.1365'
$ws.Cells.Item(136, 9).Value = 1573.3846  # I136: '1791.9' -> '1573.3846'
$ws.Cells.Item(136, 10).Value = 4151.6665  # J136: '4047.3684' -> '4151.6665'
$ws.Cells.Item(136, 11).Value = 4720.1538  # K136: '5375.700000000001' -> '4720.1538'
$ws.Cells.Item(136, 12).Value = 12454.9995  # L136: '12142.1052' -> '12454.9995'
$ws.Cells.Item(136, 13).Value = -2170.1538  # M136: '-2825.700000000001' -> '-2170.1538'
$ws.Cells.Item(136, 14).Value = -17554.9995  # N136: '-17242.1052' -> '-17554.9995'

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(119, 8).Value = 40000  # H119: '0' -> '40000'
$ws.Cells.Item(119, 10).Value = 40000  # J119: '0' -> '40000'
$ws.Cells.Item(119, 12).Value = 40000  # L119: '0' -> '40000'
$ws.Cells.Item(119, 14).Value = -49676  # N119: None -> '-49676'

$ws.Cells.Item(132, 8).Value = 1422.5416  # H132: '1574.7693' -> '1422.5416'
$ws.Cells.Item(132, 9).Value = 645.1489  # I132: '724.2381' -> '645.1489'
$ws.Cells.Item(132, 10).Value = 2884.04  # J132: '3127.913' -> '2884.04'
$ws.Cells.Item(132, 11).Value = 1935.4467  # K132: '2172.7143' -> '1935.4467'
$ws.Cells.Item(132, 12).Value = 8652.119999999999  # L132: '9383.739' -> '8652.119999999999'
$ws.Cells.Item(132, 13).Value = 594.5533  # M132: '357.2856999999999' -> '594.5533'
$ws.Cells.Item(132, 14).Value = -13712.12  # N132: '-14443.739' -> '-13712.12'

$ws.Cells.Item(136, 8).Value = 4238.3716  # H136: '5236.037' -> '4238.3716'
$ws.Cells.Item(136, 9).Value = 5082.087  # I136: '6249.8887' -> '5082.087'
$ws.Cells.Item(136, 10).Value = 2621.25  # J136: '3208.3333' -> '2621.25'
$ws.Cells.Item(136, 11).Value = 15246.261  # K136: '18749.6661' -> '15246.261'
$ws.Cells.Item(136, 12).Value = 7863.75  # L136: '9624.999899999999' -> '7863.75'
$ws.Cells.Item(136, 13).Value = -12696.261  # M136: '-16199.6661' -> '-12696.261'
$ws.Cells.Item(136, 14).Value = -12963.75  # N136: '-14724.9999' -> '-12963.75'
